$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new shared string / header for Homework 4
$ws.Range("I2").Value = "Homework 4"

# Set column I width similar to the other homework columns (~13.21 chars;
# engine stores width snapped to pixel grid, 12.33 is the closest input that
# rounds to the intended ~13.2 stored width)
$ws.Columns.Item(9).ColumnWidth = 12.33

# Fill in Homework 4 grades (column I) as formulas, rows 3-19 (row 10 has no data, like row 10 in source)
$ws.Range("I3").Formula = "=10/10"
$ws.Range("I4").Formula = "=0"
$ws.Range("I5").Formula = "=0"
$ws.Range("I6").Formula = "=11/10"
$ws.Range("I7").Formula = "=10.5/10"
$ws.Range("I8").Formula = "=10/10"
$ws.Range("I9").Formula = "=9.5/10"
$ws.Range("I11").Formula = "=10/10"
$ws.Range("I12").Formula = "=9/10"
$ws.Range("I13").Formula = "=11/10"
$ws.Range("I14").Formula = "=9/10"
$ws.Range("I15").Formula = "=8.5/10"
$ws.Range("I16").Formula = "=10/10"
$ws.Range("I17").Formula = "=9.5/10"
$ws.Range("I18").Formula = "=11/10"
$ws.Range("I19").Formula = "=10/10"

# Update selection to I10 (matches the edited workbook's saved cursor position)
$ws.Range("I10").Select() | Out-Null
